# Fix typos: pluralize "tag" -> "tags" in constraint names on the Configuration sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

$ws.Range("A12").Value = "Speaker preferred timeslot tags"
$ws.Range("A13").Value = "Speaker undesired timeslot tags"
$ws.Range("A14").Value = "Talk preferred timeslot tags"
$ws.Range("A15").Value = "Talk undesired timeslot tags"
$ws.Range("A16").Value = "Speaker preferred room tags"
$ws.Range("A17").Value = "Speaker undesired room tags"
$ws.Range("A18").Value = "Talk preferred room tags"
$ws.Range("A19").Value = "Talk undesired room tags"

$ws.Range("A29").Value = "Speaker required timeslot tags"
$ws.Range("A30").Value = "Speaker prohibited timeslot tags"
$ws.Range("A31").Value = "Talk required timeslot tags"
$ws.Range("A32").Value = "Talk prohibited timeslot tags"
$ws.Range("A33").Value = "Speaker required room tags"
$ws.Range("A34").Value = "Speaker prohibited room tags"
$ws.Range("A35").Value = "Talk required room tags"
$ws.Range("A36").Value = "Talk prohibited room tags"
$ws.Range("A37").Value = "Talk mutually-exclusive-talks tags"
